$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "ปั่นงานไม่ทันแล้วทำอย่างไรดี"
$ws.Range("B29").Value = "งาน"
$ws.Range("C29").Value = "เวลาชีวิต"

$ws.Range("A30").Value = "งานเยอะมากเลย เลือกทำไม่ถูก"
$ws.Range("B30").Value = "งาน"

$ws.Range("A31").Value = "วิธีการเลี้ยงปลาทอง"
$ws.Range("B31").Value = "ปลาทอง"
$ws.Range("C31").Value = "สัตว์เลี้ยง"
